$wb = $excel.ActiveWorkbook

# Update "Latest Handoff Datetime" (column D) for the 5e1496a8... file row (row 4)
# on the zh-cn handoff-status sheet, reflecting a new handoff generated for this report.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-17 07:27:46"

# Same update on the de-de handoff-status sheet.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-17 07:27:56"
